$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

# Percentage-formatted text cells (column H) need NumberFormat forced to
# Text ("@") before assignment, otherwise Excel auto-converts "NN%" strings
# into numeric percentage values instead of keeping literal text.

$ws.Range("E2").Value = "2026-02-15 18:18:58"
$ws.Range("E3").Value = "2026-02-15 18:19:00"
$ws.Range("E4").Value = "2026-02-15 18:19:03"
$ws.Range("O4").Value = "7.3 °C"
$ws.Range("E5").Value = "2026-02-15 18:19:06"
$ws.Range("I5").Value = "3.4 mm"
$ws.Range("O5").Value = "-4.9 °C"
$ws.Range("E6").Value = "2026-02-15 18:19:09"
$ws.Range("O6").Value = "8.5 °C"
$ws.Range("E7").Value = "2026-02-15 18:19:11"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "48%"
$ws.Range("O7").Value = "11.5 °C"
$ws.Range("E8").Value = "2026-02-15 18:19:14"
$ws.Range("E9").Value = "2026-02-15 18:19:17"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "48%"
$ws.Range("E10").Value = "2026-02-15 18:19:20"
$ws.Range("E11").Value = "2026-02-15 18:19:23"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "39%"
$ws.Range("E12").Value = "2026-02-15 18:19:25"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "52%"
$ws.Range("E13").Value = "2026-02-15 18:19:28"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "33%"
$ws.Range("E14").Value = "2026-02-15 18:19:31"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "57%"
$ws.Range("E15").Value = "2026-02-15 18:19:32"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "48%"
$ws.Range("E16").Value = "2026-02-15 18:19:35"
$ws.Range("M16").Value = "0.5 °C 17:52 TU"
$ws.Range("E17").Value = "2026-02-15 18:19:37"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "35%"
$ws.Range("E18").Value = "2026-02-15 18:19:40"
$ws.Range("O18").Value = "7.2 °C"
$ws.Range("E19").Value = "2026-02-15 18:19:43"
$ws.Range("E20").Value = "2026-02-15 18:19:45"
$ws.Range("O20").Value = "-2.9 °C"
$ws.Range("E21").Value = "2026-02-15 18:19:48"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "36%"
$ws.Range("O21").Value = "7.7 °C"
$ws.Range("E22").Value = "2026-02-15 18:19:51"
$ws.Range("N22").Value = "-6.5 °C 17:30 TU"
$ws.Range("E23").Value = "2026-02-15 18:19:54"
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = "63%"
$ws.Range("I23").Value = "1.6 mm"
$ws.Range("O23").Value = "-3.9 °C"
$ws.Range("E24").Value = "2026-02-15 18:19:57"
$ws.Range("E25").Value = "2026-02-15 18:20:00"
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H25").Value = "62%"
$ws.Range("E26").Value = "2026-02-15 18:20:03"
$ws.Range("E27").Value = "2026-02-15 18:20:05"
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H27").Value = "47%"
$ws.Range("O27").Value = "-0.2 °C"
$ws.Range("E28").Value = "2026-02-15 18:20:08"
$ws.Range("E29").Value = "2026-02-15 18:20:10"
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = "55%"
$ws.Range("E30").Value = "2026-02-15 18:20:13"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "52%"
$ws.Range("J30").Value = "1015.2 hPa"
$ws.Range("E31").Value = "2026-02-15 18:20:16"
$ws.Range("E32").Value = "2026-02-15 18:20:19"
$ws.Range("K32").Value = "9.3 MJ/m2"
$ws.Range("O32").Value = "3.4 °C"
$ws.Range("E33").Value = "2026-02-15 18:20:21"
$ws.Range("E34").Value = "2026-02-15 18:20:24"
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = "50%"
$ws.Range("O34").Value = "0.8 °C"
$ws.Range("E35").Value = "2026-02-15 18:20:27"
$ws.Range("E36").Value = "2026-02-15 18:20:30"
$ws.Range("E37").Value = "2026-02-15 18:20:33"
$ws.Range("H37").NumberFormat = "@"
$ws.Range("H37").Value = "52%"
$ws.Range("O37").Value = "6.1 °C"
$ws.Range("E38").Value = "2026-02-15 18:20:35"
$ws.Range("E39").Value = "2026-02-15 18:20:38"
$ws.Range("O39").Value = "-3.2 °C"
$ws.Range("E40").Value = "2026-02-15 18:20:41"
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = "33%"
$ws.Range("E41").Value = "2026-02-15 18:20:43"
$ws.Range("O41").Value = "12.3 °C"
$ws.Range("E42").Value = "2026-02-15 18:20:46"
$ws.Range("H42").NumberFormat = "@"
$ws.Range("H42").Value = "54%"
$ws.Range("E43").Value = "2026-02-15 18:20:49"
$ws.Range("O43").Value = "6.1 °C"
$ws.Range("E44").Value = "2026-02-15 18:20:52"
$ws.Range("H44").NumberFormat = "@"
$ws.Range("H44").Value = "76%"
$ws.Range("O44").Value = "-4.3 °C"
$ws.Range("E45").Value = "2026-02-15 18:20:55"
$ws.Range("J45").Value = "1023.6 hPa"
$ws.Range("O45").Value = "0.8 °C"
$ws.Range("E46").Value = "2026-02-15 18:20:57"
$ws.Range("O46").Value = "11.6 °C"
